$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (Volume/Number and reporting week dates)
# ---------------------------------------------------------------------------
# "Volume 30   Number  12" -> "Volume 30   Number  13"
$ws.Range("A8").Characters(22, 1).Text = "3"

# "Report Covering the Week  3/20/2023  Through  3/26/2023"
# -> "Report Covering the Week  3/27/2023  Through  4/2/2023"
$ws.Range("C9").Characters(27, 9).Text = "3/27/2023"
$ws.Range("C9").Characters(47, 9).Text = "4/2/2023"

# ---------------------------------------------------------------------------
# Weekly crime-stat table updates (rows 15-27)
# Cells that change data type/number-format (text "0"/"***.*" <-> numeric)
# are first re-based on a same-shaped template cell (which copies both the
# value and the formatting), then have their final value applied.
# ---------------------------------------------------------------------------

# Row 15 - Rape
$ws.Range("I14").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 66.666666666666
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = -16.666666666666

# Row 16 - Robbery
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 166.666666666667
$ws.Range("I16").Value = 25
$ws.Range("K16").Value = 4.166666666666
$ws.Range("L16").Value = -13.793103448275
$ws.Range("M16").Value = 19.047619047619
$ws.Range("N16").Value = -83.552631578947

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 83.333333333333
$ws.Range("I17").Value = 26
$ws.Range("J17").Value = 18
$ws.Range("K17").Value = 44.444444444444
$ws.Range("L17").Value = 85.714285714285
$ws.Range("M17").Value = 44.444444444444
$ws.Range("N17").Value = 8.333333333333

# Row 18 - Burglary
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -36.363636363636
$ws.Range("I18").Value = 39
$ws.Range("J18").Value = 37
$ws.Range("K18").Value = 5.405405405405
$ws.Range("L18").Value = 44.444444444444
$ws.Range("M18").Value = 14.705882352941
$ws.Range("N18").Value = -81.428571428571

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 60
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = -2
$ws.Range("I19").Value = 172
$ws.Range("J19").Value = 175
$ws.Range("K19").Value = -1.714285714285
$ws.Range("L19").Value = 36.507936507936
$ws.Range("M19").Value = -4.972375690607
$ws.Range("N19").Value = -66.471734892787

# Row 20 - G.L.A.
$ws.Range("I14").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 2
$ws.Range("K30").Copy($ws.Range("E20"))
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 14
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = 61.538461538461
$ws.Range("M20").Value = 320
$ws.Range("N20").Value = -93.046357615894

# Row 21 - TOTAL
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 27.777777777777
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = 6.493506493506
$ws.Range("I21").Value = 289
$ws.Range("J21").Value = 273
$ws.Range("K21").Value = 5.860805860805
$ws.Range("L21").Value = 36.320754716981
$ws.Range("M21").Value = 10.727969348659
$ws.Range("N21").Value = -76.056338028169

# Row 22 - Transit
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("G22"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("L22").Value = 16.666666666666

# Row 23 - Housing
$ws.Range("C23").Value = 1
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = 300
$ws.Range("I23").Value = 14
$ws.Range("K23").Value = 366.666666666667
$ws.Range("L23").Value = 180
$ws.Range("M23").Value = 100

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -34.782608695652
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = -17.021276595744
$ws.Range("I24").Value = 250
$ws.Range("J24").Value = 318
$ws.Range("K24").Value = -21.383647798742
$ws.Range("L24").Value = -35.233160621761
$ws.Range("M24").Value = 10.619469026548

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 22
$ws.Range("H25").Value = 57.142857142857
$ws.Range("I25").Value = 58
$ws.Range("J25").Value = 49
$ws.Range("K25").Value = 18.367346938775
$ws.Range("L25").Value = 16
$ws.Range("M25").Value = 1.754385964912

# Row 26 - UCR Rape*
$ws.Range("I14").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 6
$ws.Range("K26").Value = 20
$ws.Range("L26").Value = 20

# Row 27 - Other Sex Crimes
$ws.Range("I14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("L27").Value = 80
